# Update the PROFESORES sheet:
#  - remove the "Saúl Calderón Ramírez" (ASISTENTE) professor row
#  - duplicate the "super" administrative row
#  - add a new "Ariana B" professor row
#
# The cleanest way to reproduce the final data (and drop the old per-cell
# styling / hyperlink-blue font / custom row heights that the original
# rows carried) is to delete the whole existing data block and retype the
# final values into fresh, unformatted rows - exactly what happened in the
# real edit (compare the xml diff: none of the new rows carry any `s=`
# style attribute any more).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PROFESORES")

# Remove the six existing data rows (rows 2-7) ...
$ws.Rows("2:7").Delete()

# ... then reinsert six blank rows so the trailing formatted blank row
# (originally row 9) ends up back in the same place once seven new data
# rows have been written (rows 2-8).
$ws.Rows("3:8").Insert()

# Row 2: Ericka Solano Fernández (unchanged)
$ws.Range("A2").Value = "1-1111-1111"
$ws.Range("B2").Value = "Ericka Solano Fernández"
$ws.Range("C2").Value = "ersolano@itcr.ac.cr"
$ws.Range("D2").Value = "8111-1111"
$ws.Range("E2").Value = "masm"
$ws.Range("F2").Value = "COORDINADOR"

# Row 3: Alicia Salazar Hernández (unchanged)
$ws.Range("A3").Value = "2-2222-2222"
$ws.Range("B3").Value = "Alicia Salazar Hernández"
$ws.Range("C3").Value = "asalazar@itcr.ac.cr"
$ws.Range("D3").Value = "8222-2222"
$ws.Range("E3").Value = "dfsdf"
$ws.Range("F3").Value = "DIRECTOR"

# Row 4: Franco Quirós Ramírez (was row 5; "Saúl Calderón Ramírez" removed)
$ws.Range("A4").Value = "4-4444-4444"
$ws.Range("B4").Value = "Franco Quirós Ramírez"
$ws.Range("C4").Value = "fquiros@itcr.ac.cr"
$ws.Range("D4").Value = "8444-4444"
$ws.Range("E4").Value = "qwkqe"
$ws.Range("F4").Value = "DIRECTOR"

# Row 5: Ivannia Cerdas Quesada
$ws.Range("A5").Value = "5-5555-5555"
$ws.Range("B5").Value = "Ivannia Cerdas Quesada"
$ws.Range("C5").Value = "iquesada@itcr.ac.cr"
$ws.Range("D5").Value = "8555-5555"
$ws.Range("E5").Value = "qweqqqq"
$ws.Range("F5").Value = "DIRECTOR"

# Row 6: super / super / super (SUPERUSUARIO)
$ws.Range("A6").Value = "super"
$ws.Range("B6").Value = "super"
$ws.Range("C6").Value = "super"
$ws.Range("D6").Value = "8666-6666"
$ws.Range("E6").Value = "Disenno"
$ws.Range("F6").Value = "SUPERUSUARIO"

# Row 7: duplicate of the "super" row
$ws.Range("A7").Value = "super"
$ws.Range("B7").Value = "super"
$ws.Range("C7").Value = "super"
$ws.Range("D7").Value = "8666-6666"
$ws.Range("E7").Value = "Disenno"
$ws.Range("F7").Value = "SUPERUSUARIO"

# Row 8: new professor Ariana B
$ws.Range("A8").Value = "1-1670-0598"
$ws.Range("B8").Value = "Ariana B"
$ws.Range("C8").Value = "b@h.com"
$ws.Range("F8").Value = "DIRECTOR"

# D8 / E8 hold digit-only text ("22609987" / "123") that must stay text,
# not get auto-converted to numbers. Write them into a scratch cell that
# is formatted as Text first, then copy/paste-values them into place so
# the target cells end up as shared-string cells without picking up a
# quote-prefix style.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "22609987"
$ws.Range("Z1").Copy()
$ws.Range("D8").PasteSpecial(-4163)

$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "123"
$ws.Range("Z1").Copy()
$ws.Range("E8").PasteSpecial(-4163)

$ws.Range("Z1").Clear()
